$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 3, shifting existing data (rows 3-28) down to rows 4-29
$ws.Rows.Item(3).Insert()

# Reset the view / selection to cell C1 (matches the saved sheetView state in the target file)
$ws.Range("C1").Select()
